$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values per row (B,C,D,E,F,G)
$data = @{
    2  = @(812.25, 800.45, 802.7,  802.1,  39,  809.35)
    3  = @(905.5,  882.25, 903.55, 902.75, 72,  884.35)
    4  = @(43979,  43517.35, 43938.8, 43926.85, 24, 43708)
    5  = @(392.15, 384.85, 390.1,  389.8,  151, 386.15)
    6  = @(596.65, 589.35, 594.15, 594.3,  75,  591)
    7  = @(494.05, 484.15, 486.3,  487.25, 108, 485.35)
    8  = @(952.65, 935.8,  951.75, 950.4,  206, 943.65)
    9  = @(636,    623.1,  631,    631.85, 95,  623.9)
    10 = @(19493.9, 19415, 19480.05, 19476.9, 42, 0)
    11 = @(2345.55, 2332.15, 2335.35, 2335.65, 39, 2344.6)
    12 = @(582.35, 576.5,  582,    581.2,  303, 576.8)
    13 = @(926.35, 918.2,  925.8,  925.45, 12,  923.55)
    14 = @(653.3,  647.6,  648.3,  648.5,  100, 649)
    15 = @(120.85, 119.15, 120.2,  120.2,  375, 119.4)
    16 = @(3420.75, 3383.6, 3388.7, 3388.15, 14, 3410.5)
    17 = @(3296,   3277.25, 3289.85, 3288.2, 12, 3279)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}

$ws.Range("N9").Select()
